$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 "25.976.89"
Set-TextValue 2 5 "  -1.39%  "

Set-TextValue 3 4 "1.640.33"
Set-TextValue 3 5 "  -1.54%  "

Set-TextValue 4 4 "1.006"

Set-TextValue 5 4 "215.50"
Set-TextValue 5 5 "  -1.43%  "

Set-TextValue 6 4 "0.5045"
Set-TextValue 6 5 "  -2.29%  "

Set-TextValue 7 4 "1.006"
Set-TextValue 7 5 "  -0.14%  "

Set-TextValue 8 4 "0.2579"
Set-TextValue 8 5 "  +0.56%  "

Set-TextValue 9 4 "0.06441"
Set-TextValue 9 5 "  -0.15%  "

Set-TextValue 10 4 "19.55"
Set-TextValue 10 5 "  -1.89%  "

Set-TextValue 11 4 "0.07730"
Set-TextValue 11 5 "  +0.81%  "

$ws.Cells.Item(12, 2).Value = "Polkadot"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue 12 4 "4.250"
Set-TextValue 12 5 "  -1.47%  "

$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue 13 4 "1.637.54"
Set-TextValue 13 5 "  -1.81%  "

Set-TextValue 14 4 "1.862.92"
Set-TextValue 14 5 "  -1.72%  "

Set-TextValue 15 4 "0.5462"
Set-TextValue 15 5 "  -1.45%  "

Set-TextValue 16 4 "0.0₅7962"
Set-TextValue 16 5 "  -0.99%  "

Set-TextValue 17 4 "63.64"
Set-TextValue 17 5 "  -1.26%  "

Set-TextValue 18 4 "25.968.73"
Set-TextValue 18 5 "  -1.57%  "

Set-TextValue 19 4 "1.006"

Set-TextValue 20 4 "205.14"
Set-TextValue 20 5 "  -2.47%  "

Set-TextValue 21 4 "4.318"
Set-TextValue 21 5 "  -1.91%  "

Set-TextValue 22 4 "10.000"
Set-TextValue 22 5 "  -1.08%  "

Set-TextValue 23 4 "5.960"
Set-TextValue 23 5 "  +1.04%  "

Set-TextValue 24 5 "  -0.13%  "

Set-TextValue 25 4 "1.915"
Set-TextValue 25 5 "  +9.21%  "

Set-TextValue 26 4 "141.15"
Set-TextValue 26 5 "  -2.30%  "

Set-TextValue 27 4 "0.1157"
Set-TextValue 27 5 "  -0.38%  "

Set-TextValue 28 4 "15.83"
Set-TextValue 28 5 "  +0.38%  "

Set-TextValue 29 4 "6.764"
Set-TextValue 29 5 "  -3.29%  "

Set-TextValue 30 4 "0.05065"
Set-TextValue 30 5 "  -3.58%  "

Set-TextValue 31 4 "1.239"
Set-TextValue 31 5 "  -1.80%  "

Set-TextValue 32 4 "3.269"
Set-TextValue 32 5 "  -3.30%  "

Set-TextValue 33 4 "3.195"
Set-TextValue 33 5 "  -0.80%  "

Set-TextValue 34 4 "1.540"
Set-TextValue 34 5 "  -1.60%  "

Set-TextValue 35 4 "2.336"
Set-TextValue 35 5 "  -1.72%  "

Set-TextValue 36 4 "0.8964"

Set-TextValue 37 4 "2.619"
Set-TextValue 37 5 "  -5.02%  "

Set-TextValue 38 4 "0.5660"
Set-TextValue 38 5 "  -1.19%  "

Set-TextValue 39 4 "1.140.74"
Set-TextValue 39 5 "  -1.23%  "

Set-TextValue 40 4 "0.01565"
Set-TextValue 40 5 "  -2.01%  "

Set-TextValue 41 4 "2.556"
Set-TextValue 41 5 "  -1.05%  "

Set-TextValue 42 5 "  -0.16%  "

Set-TextValue 43 4 "5.640"
Set-TextValue 43 5 "  -0.19%  "

Set-TextValue 44 4 "0.8180"
Set-TextValue 44 5 "  -3.30%  "

Set-TextValue 45 4 "99.43"
Set-TextValue 45 5 "  -0.45%  "

Set-TextValue 46 4 "1.774.35"
Set-TextValue 46 5 "  -1.73%  "

Set-TextValue 47 5 "  -0.83%  "

Set-TextValue 48 4 "0.4522"
Set-TextValue 48 5 "  +0.44%  "

Set-TextValue 49 4 "1.010"
Set-TextValue 49 5 "  +0.35%  "

Set-TextValue 50 4 "54.87"
Set-TextValue 50 5 "  -2.22%  "

Set-TextValue 51 4 "0.05032"
Set-TextValue 51 5 "  -1.33%  "
